$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.700661
$ws.Range("H2").Value = 68.10198299999999
$ws.Range("I2").Value = 0.08615268874617349
$ws.Range("J2").Value = 0.08615268874617349
$ws.Range("M2").Value = 576.300578
$ws.Range("N2").Value = 1728.901734
$ws.Range("O2").Value = 0.8614732012478776
$ws.Range("P2").Value = 0.8614732012478775
$ws.Range("Q2").Value = 13082.40405528206
$ws.Range("R2").Value = 117741.6364975385
$ws.Range("S2").Value = 0.07421823257027807
$ws.Range("T2").Value = 0.07421823257027807
$ws.Range("G3").Value = 22.700661
$ws.Range("H3").Value = 68.10198299999999
$ws.Range("I3").Value = 0.08615268874617349
$ws.Range("J3").Value = 0.08615268874617349
$ws.Range("O3").Value = 0.001786049553652741
$ws.Range("P3").Value = 0.001786049553652741
$ws.Range("Q3").Value = 27.123097839602
$ws.Range("R3").Value = 244.1078805564179
$ws.Range("S3").Value = 0.0001538729712810867
$ws.Range("T3").Value = 0.0001538729712810867
$ws.Range("G4").Value = 22.700661
$ws.Range("H4").Value = 68.10198299999999
$ws.Range("I4").Value = 0.08615268874617349
$ws.Range("J4").Value = 0.08615268874617349
$ws.Range("M4").Value = 34.99993866666667
$ws.Range("N4").Value = 104.999816
$ws.Range("O4").Value = 0.05231906813505349
$ws.Range("P4").Value = 0.05231906813505348
$ws.Range("Q4").Value = 794.5217426927918
$ws.Range("R4").Value = 7150.695684235126
$ws.Range("S4").Value = 0.004507428392529106
$ws.Range("T4").Value = 0.004507428392529106
$ws.Range("G5").Value = 22.700661
$ws.Range("H5").Value = 68.10198299999999
$ws.Range("I5").Value = 0.08615268874617349
$ws.Range("J5").Value = 0.08615268874617349
$ws.Range("M5").Value = 56.47565533333333
$ws.Range("N5").Value = 169.426966
$ws.Range("O5").Value = 0.08442168106341624
$ws.Range("P5").Value = 0.08442168106341623
$ws.Range("Q5").Value = 1282.034706474842
$ws.Range("R5").Value = 11538.31235827358
$ws.Range("S5").Value = 0.007273154812085228
$ws.Range("T5").Value = 0.007273154812085227
$ws.Range("I6").Value = 0.5030288587986086
$ws.Range("J6").Value = 0.5030288587986087
$ws.Range("M6").Value = 576.300578
$ws.Range("N6").Value = 1728.901734
$ws.Range("O6").Value = 0.8614732012478776
$ws.Range("P6").Value = 0.8614732012478775
$ws.Range("Q6").Value = 76385.62275937223
$ws.Range("R6").Value = 687470.6048343502
$ws.Range("S6").Value = 0.433345881309304
$ws.Range("T6").Value = 0.433345881309304
$ws.Range("I7").Value = 0.5030288587986086
$ws.Range("J7").Value = 0.5030288587986087
$ws.Range("O7").Value = 0.001786049553652741
$ws.Range("P7").Value = 0.001786049553652741
$ws.Range("S7").Value = 0.0008984344687317026
$ws.Range("T7").Value = 0.0008984344687317026
$ws.Range("I8").Value = 0.5030288587986086
$ws.Range("J8").Value = 0.5030288587986087
$ws.Range("M8").Value = 34.99993866666667
$ws.Range("N8").Value = 104.999816
$ws.Range("O8").Value = 0.05231906813505349
$ws.Range("P8").Value = 0.05231906813505348
$ws.Range("Q8").Value = 4639.058528921284
$ws.Range("R8").Value = 41751.52676029156
$ws.Range("S8").Value = 0.0263180011373826
$ws.Range("T8").Value = 0.02631800113738261
$ws.Range("I9").Value = 0.5030288587986086
$ws.Range("J9").Value = 0.5030288587986087
$ws.Range("M9").Value = 56.47565533333333
$ws.Range("N9").Value = 169.426966
$ws.Range("O9").Value = 0.08442168106341624
$ws.Range("P9").Value = 0.08442168106341623
$ws.Range("Q9").Value = 7485.552276125478
$ws.Range("R9").Value = 67369.97048512931
$ws.Range("S9").Value = 0.04246654188319038
$ws.Range("T9").Value = 0.04246654188319038
$ws.Range("G10").Value = 41.94534433333333
$ws.Range("H10").Value = 125.836033
$ws.Range("I10").Value = 0.159189381961201
$ws.Range("J10").Value = 0.159189381961201
$ws.Range("M10").Value = 576.300578
$ws.Range("N10").Value = 1728.901734
$ws.Range("O10").Value = 0.8614732012478776
$ws.Range("P10").Value = 0.8614732012478775
$ws.Range("Q10").Value = 24173.12618370902
$ws.Range("R10").Value = 217558.1356533812
$ws.Range("S10").Value = 0.137137386482787
$ws.Range("T10").Value = 0.1371373864827869
$ws.Range("G11").Value = 41.94534433333333
$ws.Range("H11").Value = 125.836033
$ws.Range("I11").Value = 0.159189381961201
$ws.Range("J11").Value = 0.159189381961201
$ws.Range("O11").Value = 0.001786049553652741
$ws.Range("P11").Value = 0.001786049553652741
$ws.Range("Q11").Value = 50.11694057141311
$ws.Range("R11").Value = 451.0524651427179
$ws.Range("S11").Value = 0.0002843201245980588
$ws.Range("T11").Value = 0.0002843201245980587
$ws.Range("G12").Value = 41.94534433333333
$ws.Range("H12").Value = 125.836033
$ws.Range("I12").Value = 0.159189381961201
$ws.Range("J12").Value = 0.159189381961201
$ws.Range("M12").Value = 34.99993866666667
$ws.Range("N12").Value = 104.999816
$ws.Range("O12").Value = 0.05231906813505349
$ws.Range("P12").Value = 0.05231906813505348
$ws.Range("Q12").Value = 1468.084479018881
$ws.Range("R12").Value = 13212.76031116993
$ws.Range("S12").Value = 0.008328640121205129
$ws.Range("T12").Value = 0.008328640121205128
$ws.Range("G13").Value = 41.94534433333333
$ws.Range("H13").Value = 125.836033
$ws.Range("I13").Value = 0.159189381961201
$ws.Range("J13").Value = 0.159189381961201
$ws.Range("M13").Value = 56.47565533333333
$ws.Range("N13").Value = 169.426966
$ws.Range("O13").Value = 0.08442168106341624
$ws.Range("P13").Value = 0.08442168106341623
$ws.Range("Q13").Value = 2368.890809407319
$ws.Range("R13").Value = 21320.01728466588
$ws.Range("S13").Value = 0.01343903523261086
$ws.Range("T13").Value = 0.01343903523261086
$ws.Range("G14").Value = 66.302588
$ws.Range("H14").Value = 198.907764
$ws.Range("I14").Value = 0.2516290704940168
$ws.Range("J14").Value = 0.2516290704940168
$ws.Range("M14").Value = 576.300578
$ws.Range("N14").Value = 1728.901734
$ws.Range("O14").Value = 0.8614732012478776
$ws.Range("P14").Value = 0.8614732012478775
$ws.Range("Q14").Value = 38210.21978729586
$ws.Range("R14").Value = 343891.9780856628
$ws.Range("S14").Value = 0.2167717008855086
$ws.Range("T14").Value = 0.2167717008855085
$ws.Range("G15").Value = 66.302588
$ws.Range("H15").Value = 198.907764
$ws.Range("I15").Value = 0.2516290704940168
$ws.Range("J15").Value = 0.2516290704940168
$ws.Range("O15").Value = 0.001786049553652741
$ws.Range("P15").Value = 0.001786049553652741
$ws.Range("Q15").Value = 79.21934878208266
$ws.Range("R15").Value = 712.9741390387439
$ws.Range("S15").Value = 0.0004494219890418929
$ws.Range("T15").Value = 0.0004494219890418928
$ws.Range("G16").Value = 66.302588
$ws.Range("H16").Value = 198.907764
$ws.Range("I16").Value = 0.2516290704940168
$ws.Range("J16").Value = 0.2516290704940168
$ws.Range("M16").Value = 34.99993866666667
$ws.Range("N16").Value = 104.999816
$ws.Range("O16").Value = 0.05231906813505349
$ws.Range("P16").Value = 0.05231906813505348
$ws.Range("Q16").Value = 2320.586513441269
$ws.Range("R16").Value = 20885.27862097142
$ws.Range("S16").Value = 0.01316499848393664
$ws.Range("T16").Value = 0.01316499848393664
$ws.Range("G17").Value = 66.302588
$ws.Range("H17").Value = 198.907764
$ws.Range("I17").Value = 0.2516290704940168
$ws.Range("J17").Value = 0.2516290704940168
$ws.Range("M17").Value = 56.47565533333333
$ws.Range("N17").Value = 169.426966
$ws.Range("O17").Value = 0.08442168106341624
$ws.Range("P17").Value = 0.08442168106341623
$ws.Range("Q17").Value = 3744.482107596002
$ws.Range("R17").Value = 33700.33896836402
$ws.Range("S17").Value = 0.02124294913552977
$ws.Range("T17").Value = 0.02124294913552977
